$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template's two header cells had their shared-string text swapped:
#   A1 was "Item Category Name" -> now "ID"
#   B1 was "Item Group Name"    -> now "Item Category Name"
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Item Category Name"

# The header row's font is now bold (fill/color stay the same).
$ws.Range("A1:B1").Font.Bold = $true

# Column widths were adjusted: column A narrower, column B wider
# (COM ColumnWidth is quantized to the sheet's character grid, so these
# inputs are chosen to land as close as possible to the saved widths of
# 43.140625 and 47).
$ws.Columns(1).ColumnWidth = 42.333333333333336
$ws.Columns(2).ColumnWidth = 46.166666666666664
